$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142, shifting existing rows 142-144 down to 143-145.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new weekly record.
$ws.Range("A142").Value = 5
$ws.Range("B142").Value = "Macroferia Regional de Talca"
$ws.Range("C142").Value = "Maule"
$ws.Range("D142").Value = 45239
$ws.Range("D142").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E142").Value = 7
$ws.Range("F142").Value = 100112026
$ws.Range("G142").Value = "Haba"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 500
$ws.Range("K142").Value = 9000
$ws.Range("L142").Value = 9000
$ws.Range("M142").Value = 9000
$ws.Range("N142").Value = "`$/saco 25 kilos"
$ws.Range("O142").Value = "Región del Maule"
$ws.Range("P142").Value = 360
$ws.Range("Q142").Value = 25
$ws.Range("R142").Value = "Hortaliza"

# Row 143 (previously row 142): update date and price values, and origin region.
$ws.Range("D143").Value = 44858
$ws.Range("K143").Value = 7000
$ws.Range("L143").Value = 7000
$ws.Range("M143").Value = 7000
$ws.Range("O143").Value = "Región del Maule"
$ws.Range("P143").Value = 280

# Row 144 (previously row 143): update date, price values, origin region stays O'Higgins.
$ws.Range("D144").Value = 44477
$ws.Range("K144").Value = 8000
$ws.Range("L144").Value = 8000
$ws.Range("M144").Value = 8000
$ws.Range("O144").Value = "Región de O'Higgins"
$ws.Range("P144").Value = 320

# Row 145 (previously row 144): update date; price values/origin remain the same.
$ws.Range("D145").Value = 44508
$ws.Range("K145").Value = 6000
$ws.Range("L145").Value = 6000
$ws.Range("M145").Value = 6000
$ws.Range("O145").Value = "Región del Maule"
$ws.Range("P145").Value = 240

"ok"
